$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B2 from a numeric value to a text value "154"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "154"

# Add a new row with Jorge's data
$ws.Range("A3").Value = "Jorge"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "158"

# C3 holds the same date as C2 (2001-08-17) but formatted as date only (no time)
$ws.Range("C3").Value = (Get-Date -Year 2001 -Month 8 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C3").NumberFormat = "YYYY-MM-DD"
